# Apply "想去人数" (want-to-go count) updates to the 展览 (sheet 1)
# and 全部类型 (sheet 4) worksheets, per the commit's generated-data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1 (展览) ---
$ws1.Range("F4").Value  = 1241
$ws1.Range("F10").Value = 3453
$ws1.Range("F11").Value = 126
$ws1.Range("F13").Value = 68
$ws1.Range("F24").Value = 2590
$ws1.Range("F25").Value = 5104
$ws1.Range("F27").Value = 73
$ws1.Range("F29").Value = 1301
$ws1.Range("F35").Value = 112
$ws1.Range("F39").Value = 794
$ws1.Range("F43").Value = 474

# --- Sheet 4 (全部类型) ---
$ws4.Range("F4").Value  = 1241
$ws4.Range("F10").Value = 3453
$ws4.Range("F11").Value = 126
$ws4.Range("F13").Value = 68
$ws4.Range("F25").Value = 2590
$ws4.Range("F26").Value = 5104
$ws4.Range("F28").Value = 73
$ws4.Range("F30").Value = 1301
$ws4.Range("F36").Value = 112
$ws4.Range("F40").Value = 794
$ws4.Range("F44").Value = 474
